$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) column values
$ws.Range('D2').Value = '26.041.08'
$ws.Range('D3').Value = '1.649.87'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '218.30'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5204'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2635'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06323'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.38'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07670'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.581'
$ws.Range('D13').Value = '1.650.78'
$ws.Range('D14').Value = '1.877.12'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5598'
$ws.Range('D16').Value = '0.0₅8129'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '65.19'
$ws.Range('D18').Value = '26.033.60'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.004'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.620'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.48'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '191.63'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.900'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '143.85'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1186'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.196'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.87'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.512'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05422'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.448'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.346'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.554'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.425'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5626'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01578'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.868'
$ws.Range('D42').Value = '1.028.00'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8254'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '100.95'
$ws.Range('D45').Value = '1.785.52'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '57.37'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.9988'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.4324'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.938'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05140'

# Update Volume(1h) (E) column values
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('E11').Value = '  -1.45%  '
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('E21').Value = '  +4.16%  '
$ws.Range('E22').Value = '  -0.96%  '
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('E24').Value = '  -0.28%  '
$ws.Range('E25').Value = '  -1.92%  '
$ws.Range('E26').Value = '  -1.57%  '
$ws.Range('E27').Value = '  +0.52%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('E30').Value = '  -3.49%  '
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('E32').Value = '  -0.94%  '
$ws.Range('E34').Value = '  -2.54%  '
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('E42').Value = '  -3.15%  '
$ws.Range('E43').Value = '  -1.59%  '
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('E46').Value = '  +5.64%  '
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('E51').Value = '  -3.39%  '
